$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7733245491981506
$ws.Range("B1").Value = 1.45144510269165
$ws.Range("C1").Value = 5.572677612304688
$ws.Range("D1").Value = 3.164451599121094
$ws.Range("E1").Value = 1.497875452041626
